# pushOfbiz credentials update: the ofbiz user becomes "abcd" and the
# ofbiz password cell becomes the number 1234 (its original text,
# "C@bi$ush5", is kept only as the displayed text of the existing
# mailto hyperlink on that cell). Selection cursor moves to E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldOfbizPass = "C@bi`$ush5"

# --- Keep F2's hyperlink (it still points at "mailto:C@bi$ush5") but give
# it an explicit display text, since the cell itself is about to stop
# holding that string as its literal value/text.
$f2Range = $ws.Range("F2")
foreach ($h in $f2Range.Hyperlinks) {
    $h.TextToDisplay = $oldOfbizPass
}
foreach ($h in $ws.Hyperlinks) {
    if ($h.Address -eq "") {
        $h.Address = "mailto:$oldOfbizPass"
    }
}
$oldF2Links = @()
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$F$2') {
        $oldF2Links += $h
    }
}
if ($oldF2Links.Count -gt 1) {
    $oldF2Links[0].Delete()
}

# --- Actual credential edits ---
$ws.Range("E2").Value = "abcd"
$ws.Range("F2").Value = 1234

# --- Cursor/selection moves from J2 to E2 ---
$ws.Range("E2").Select()
